# Auto-generated edit script applying crypto price/volume updates
# and two coin-row swaps (39<->40, 49<->50), per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing Text format so Excel does not
# auto-coerce numeric-looking strings (e.g. "317.32") into Number type,
# matching the original inlineStr/text storage of these cells.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value2 = $value
}

# --- Price (D) / Volume(1h) (E) updates ---
Set-TextValue $ws.Range("D2") "42.901.87"
Set-TextValue $ws.Range("E2") "  -4.76%  "
Set-TextValue $ws.Range("D3") "2.222.49"
Set-TextValue $ws.Range("E3") "  -5.80%  "
Set-TextValue $ws.Range("E4") "  -0.14%  "
Set-TextValue $ws.Range("D5") "317.32"
Set-TextValue $ws.Range("E5") "  +1.82%  "
Set-TextValue $ws.Range("D6") "99.75"
Set-TextValue $ws.Range("E6") "  -7.39%  "
Set-TextValue $ws.Range("D7") "0.593"
Set-TextValue $ws.Range("E7") "  -5.94%  "
Set-TextValue $ws.Range("D8") "0.999"
Set-TextValue $ws.Range("E8") "  -0.20%  "
Set-TextValue $ws.Range("D9") "0.564"
Set-TextValue $ws.Range("E9") "  -7.10%  "
Set-TextValue $ws.Range("D10") "37.27"
Set-TextValue $ws.Range("E10") "  -8.49%  "
Set-TextValue $ws.Range("D11") "54.07"
Set-TextValue $ws.Range("E11") "  -2.63%  "
Set-TextValue $ws.Range("D12") "0.0831"
Set-TextValue $ws.Range("E12") "  -8.97%  "
Set-TextValue $ws.Range("D13") "7.82"
Set-TextValue $ws.Range("E13") "  -6.97%  "
Set-TextValue $ws.Range("E14") "  -2.86%  "
Set-TextValue $ws.Range("D15") "0.865"
Set-TextValue $ws.Range("E15") "  -10.77%  "
Set-TextValue $ws.Range("D16") "2.554.29"
Set-TextValue $ws.Range("E16") "  -6.06%  "
Set-TextValue $ws.Range("D17") "14.29"
Set-TextValue $ws.Range("E17") "  -5.84%  "
Set-TextValue $ws.Range("D18") "2.229.32"
Set-TextValue $ws.Range("E18") "  -5.48%  "
Set-TextValue $ws.Range("D19") "42.758.35"
Set-TextValue $ws.Range("E19") "  -4.87%  "
Set-TextValue $ws.Range("D20") "15.11"
Set-TextValue $ws.Range("E20") "  +6.53%  "
Set-TextValue $ws.Range("D21") "0.0₃0967"
Set-TextValue $ws.Range("E21") "  -8.50%  "
Set-TextValue $ws.Range("D22") "6.47"
Set-TextValue $ws.Range("E22") "  -9.98%  "
Set-TextValue $ws.Range("D23") "65.50"
Set-TextValue $ws.Range("E23") "  -10.09%  "
Set-TextValue $ws.Range("D24") "3.16"
Set-TextValue $ws.Range("E24") "  -9.77%  "
Set-TextValue $ws.Range("D25") "236.44"
Set-TextValue $ws.Range("E25") "  -8.31%  "
Set-TextValue $ws.Range("D26") "2.14"
Set-TextValue $ws.Range("E26") "  -6.80%  "
Set-TextValue $ws.Range("D27") "0.999"
Set-TextValue $ws.Range("E27") "  -0.24%  "
Set-TextValue $ws.Range("D28") "10.12"
Set-TextValue $ws.Range("E28") "  -8.50%  "
Set-TextValue $ws.Range("E29") "  -4.85%  "
Set-TextValue $ws.Range("D30") "6.42"
Set-TextValue $ws.Range("E30") "  -10.12%  "
Set-TextValue $ws.Range("D31") "0.0903"
Set-TextValue $ws.Range("E31") "  -6.23%  "
Set-TextValue $ws.Range("D32") "20.52"
Set-TextValue $ws.Range("E32") "  -7.79%  "
Set-TextValue $ws.Range("D33") "34.30"
Set-TextValue $ws.Range("E33") "  -7.23%  "
Set-TextValue $ws.Range("D34") "156.70"
Set-TextValue $ws.Range("D35") "2.78"
Set-TextValue $ws.Range("E35") "  -6.19%  "
Set-TextValue $ws.Range("D36") "3.22"
Set-TextValue $ws.Range("E36") "  +11.29%  "
Set-TextValue $ws.Range("D37") "1.97"
Set-TextValue $ws.Range("E37") "  +13.79%  "
Set-TextValue $ws.Range("E38") "  -5.68%  "
Set-TextValue $ws.Range("D41") "0.105"
Set-TextValue $ws.Range("E41") "  -9.41%  "
Set-TextValue $ws.Range("D42") "0.0326"
Set-TextValue $ws.Range("E42") "  -7.27%  "
Set-TextValue $ws.Range("D43") "1.941.36"
Set-TextValue $ws.Range("E43") "  +3.38%  "
Set-TextValue $ws.Range("E44") "  +0.08%  "
Set-TextValue $ws.Range("D45") "12.47"
Set-TextValue $ws.Range("E45") "  -2.15%  "
Set-TextValue $ws.Range("D46") "88.67"
Set-TextValue $ws.Range("E46") "  -10.93%  "
Set-TextValue $ws.Range("D47") "0.208"
Set-TextValue $ws.Range("E47") "  -8.53%  "
Set-TextValue $ws.Range("D48") "5.39"
Set-TextValue $ws.Range("E48") "  -3.95%  "
Set-TextValue $ws.Range("D51") "0.878"
Set-TextValue $ws.Range("E51") "  +19.39%  "

# --- Rows whose coin swapped rank position (full row content changes) ---
Set-TextValue $ws.Range("B39") "NEARProtocol"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D39") "3.95"
Set-TextValue $ws.Range("E39") "  +0.85%  "
Set-TextValue $ws.Range("B40") "RenderToken"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D40") "4.46"
Set-TextValue $ws.Range("E40") "  -4.68%  "
Set-TextValue $ws.Range("B49") "ordi"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-TextValue $ws.Range("D49") "76.19"
Set-TextValue $ws.Range("E49") "  -5.68%  "
Set-TextValue $ws.Range("B50") "MultiversX"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue $ws.Range("D50") "60.54"
Set-TextValue $ws.Range("E50") "  -12.69%  "
